# Issue #4 Image Browser should show directory
#
# Mark Issue #4 (row 5) and Issue #5 (row 6) as "DONE" in the Status
# column (B) of the "Issues" sheet, update the current selection, and
# nudge the saved window position to reflect the author's session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Status column (B) gets a new "DONE" value for issue #4 and issue #5.
$ws.Range("B5").Value = "DONE"
$ws.Range("B6").Value = "DONE"

# Row 6 now wraps onto a second line once the Status column is
# populated, so its height grows to match row 5's wrapped height.
$ws.Rows.Item(6).RowHeight = 29

# The active cell moves to C5 as part of reviewing the newly
# completed issue.
$ws.Range("C5").Select()

# The workbook window was repositioned on screen when last saved.
$excel.ActiveWindow.Left = 2140
